# Apply stimulus updates:
#  1. Rename face stimuli to book stimuli ("face//face_NN.jpg" -> "book//book_NN.jpg")
#     across the prompt/correct/distractor file columns (A-D).
#  2. Expand abbreviated correct_ans codes in column L to full words
#     (r -> right, y -> left, b -> center).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1
$firstRow = $used.Row

$ansMap = @{ "r" = "right"; "y" = "left"; "b" = "center" }

for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    # Columns A, B, C, D hold image file references that may need the
    # face -> book rename.
    for ($c = 1; $c -le 4; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null -and $val -is [string] -and $val.Contains("face//face_")) {
            $cell.Value2 = $val.Replace("face//face_", "book//book_")
        }
    }

    # Column L (12) holds the abbreviated correct answer code.
    $lCell = $ws.Cells.Item($r, 12)
    $lVal = $lCell.Value2
    if ($lVal -ne $null -and $ansMap.ContainsKey([string]$lVal)) {
        $lCell.Value2 = $ansMap[[string]$lVal]
    }
}
